$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(105, 8).Value = 24835.5
$ws.Cells.Item(105, 10).Value = 24835.5
$ws.Cells.Item(105, 12).Value = 24835.5
$ws.Cells.Item(105, 14).Value = -31823.5

$ws.Cells.Item(138, 8).Value = 2380.449
$ws.Cells.Item(138, 9).Value = 1342.4762
$ws.Cells.Item(138, 10).Value = 3158.9285
$ws.Cells.Item(138, 11).Value = 4027.4286
$ws.Cells.Item(138, 12).Value = 9476.7855
$ws.Cells.Item(138, 13).Value = 1112.5714
$ws.Cells.Item(138, 14).Value = -19756.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 21866.615
$ws.Cells.Item(32, 9).Value = 4991.1875
$ws.Cells.Item(32, 11).Value = 4991.1875
$ws.Cells.Item(32, 13).Value = -4704.1875

$ws.Cells.Item(51, 8).Value = 61735
$ws.Cells.Item(51, 10).Value = 61735
$ws.Cells.Item(51, 12).Value = 61735
$ws.Cells.Item(51, 14).Value = -63247

$ws.Cells.Item(110, 8).Value = 873.6667
$ws.Cells.Item(110, 9).Value = 877.4706
$ws.Cells.Item(110, 10).Value = 857.5
$ws.Cells.Item(110, 11).Value = 877.4706
$ws.Cells.Item(110, 12).Value = 857.5
$ws.Cells.Item(110, 13).Value = 1167.5294
$ws.Cells.Item(110, 14).Value = -4947.5

$ws.Cells.Item(122, 8).Value = 1488.8125
$ws.Cells.Item(122, 9).Value = 935.6667
$ws.Cells.Item(122, 10).Value = 2200
$ws.Cells.Item(122, 11).Value = 2807.0001
$ws.Cells.Item(122, 12).Value = 6600
$ws.Cells.Item(122, 13).Value = -357.0001000000002
$ws.Cells.Item(122, 14).Value = -11500

$ws.Cells.Item(132, 8).Value = 2463.653
$ws.Cells.Item(132, 9).Value = 2077.1538
$ws.Cells.Item(132, 10).Value = 3971
$ws.Cells.Item(132, 11).Value = 6231.4614
$ws.Cells.Item(132, 12).Value = 11913
$ws.Cells.Item(132, 13).Value = -3701.4614
$ws.Cells.Item(132, 14).Value = -16973

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3869.318
$ws.Cells.Item(86, 9).Value = 3377.9412
$ws.Cells.Item(86, 10).Value = 5540
$ws.Cells.Item(86, 11).Value = 3377.9412
$ws.Cells.Item(86, 12).Value = 5540
$ws.Cells.Item(86, 13).Value = -2254.9412
$ws.Cells.Item(86, 14).Value = -7786

$ws.Cells.Item(89, 8).Value = 3869.318
$ws.Cells.Item(89, 9).Value = 3377.9412
$ws.Cells.Item(89, 10).Value = 5540
$ws.Cells.Item(89, 11).Value = 16889.706
$ws.Cells.Item(89, 12).Value = 27700
$ws.Cells.Item(89, 13).Value = -11273.706
$ws.Cells.Item(89, 14).Value = -38932

$ws.Cells.Item(107, 8).Value = 1305.9166
$ws.Cells.Item(107, 9).Value = 1377.1
$ws.Cells.Item(107, 10).Value = 950
$ws.Cells.Item(107, 11).Value = 1377.1
$ws.Cells.Item(107, 12).Value = 950
$ws.Cells.Item(107, 13).Value = 542.9000000000001
$ws.Cells.Item(107, 14).Value = -4790

$ws.Cells.Item(134, 8).Value = 2278.7144
$ws.Cells.Item(134, 9).Value = 1888.8
$ws.Cells.Item(134, 10).Value = 3253.5
$ws.Cells.Item(134, 11).Value = 5666.4
$ws.Cells.Item(134, 12).Value = 9760.5
$ws.Cells.Item(134, 13).Value = -3131.4
$ws.Cells.Item(134, 14).Value = -14830.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1000
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 1000
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 1000
$ws.Cells.Item(16, 13).ClearContents()
$ws.Cells.Item(16, 14).Value = -1574

$ws.Cells.Item(31, 8).Value = 1862.641
$ws.Cells.Item(31, 9).Value = 1399.6207
$ws.Cells.Item(31, 10).Value = 3205.4
$ws.Cells.Item(31, 11).Value = 1399.6207
$ws.Cells.Item(31, 12).Value = 3205.4
$ws.Cells.Item(31, 13).Value = -1104.6207
$ws.Cells.Item(31, 14).Value = -3795.4

$ws.Cells.Item(34, 8).Value = 1862.641
$ws.Cells.Item(34, 9).Value = 1399.6207
$ws.Cells.Item(34, 10).Value = 3205.4
$ws.Cells.Item(34, 11).Value = 1399.6207
$ws.Cells.Item(34, 12).Value = 3205.4
$ws.Cells.Item(34, 13).Value = -1197.6207
$ws.Cells.Item(34, 14).Value = -3609.4

$ws.Cells.Item(113, 8).Value = 1000
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 1000
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 1000
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1799.5454
$ws.Cells.Item(113, 9).Value = 5411.6665
$ws.Cells.Item(113, 10).Value = 445
$ws.Cells.Item(113, 11).Value = 16234.9995
$ws.Cells.Item(113, 12).Value = 1335
$ws.Cells.Item(113, 13).Value = -14064.9995
$ws.Cells.Item(113, 14).Value = -5675

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(105, 8).Value = 19333.334
$ws.Cells.Item(105, 10).Value = 19333.334
$ws.Cells.Item(105, 12).Value = 19333.334
$ws.Cells.Item(105, 14).Value = -26321.334

$ws.Cells.Item(126, 8).Value = 1442.2
$ws.Cells.Item(126, 9).Value = 655.5
$ws.Cells.Item(126, 11).Value = 1966.5
$ws.Cells.Item(126, 13).Value = 503.5

$ws.Cells.Item(132, 8).Value = 2716.4119
$ws.Cells.Item(132, 9).Value = 2185.9473
$ws.Cells.Item(132, 10).Value = 3388.3333
$ws.Cells.Item(132, 11).Value = 6557.841899999999
$ws.Cells.Item(132, 12).Value = 10164.9999
$ws.Cells.Item(132, 13).Value = -4027.841899999999
$ws.Cells.Item(132, 14).Value = -15224.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(48, 8).Value = 22500
$ws.Cells.Item(48, 10).Value = 40000
$ws.Cells.Item(48, 12).Value = 40000
$ws.Cells.Item(48, 14).Value = -41322

$ws.Cells.Item(55, 8).Value = 243.52
$ws.Cells.Item(55, 9).Value = 199.05882
$ws.Cells.Item(55, 10).Value = 338
$ws.Cells.Item(55, 11).Value = 199.05882
$ws.Cells.Item(55, 12).Value = 338
$ws.Cells.Item(55, 13).Value = -26.05882
$ws.Cells.Item(55, 14).Value = -684

$ws.Cells.Item(122, 8).Value = 3172.2258
$ws.Cells.Item(122, 9).Value = 3123.2173
$ws.Cells.Item(122, 11).Value = 9369.651899999999
$ws.Cells.Item(122, 13).Value = -6919.651899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(7, 8).Value = 60000
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 60000
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 60000
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(7, 14).Value = -60226

$ws.Cells.Item(8, 8).Value = 7151391.5
$ws.Cells.Item(8, 9).Value = 12501363
$ws.Cells.Item(8, 10).Value = 18096
$ws.Cells.Item(8, 11).Value = 12501363
$ws.Cells.Item(8, 12).Value = 18096
$ws.Cells.Item(8, 13).Value = -12501223
$ws.Cells.Item(8, 14).Value = -18376

$ws.Cells.Item(9, 8).Value = 4000
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 4000
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 4000
$ws.Cells.Item(9, 13).ClearContents()
$ws.Cells.Item(9, 14).Value = -4280

$ws.Cells.Item(11, 8).Value = 21000
$ws.Cells.Item(11, 9).Value = 15000
$ws.Cells.Item(11, 10).Value = 24000
$ws.Cells.Item(11, 11).Value = 15000
$ws.Cells.Item(11, 12).Value = 24000
$ws.Cells.Item(11, 13).Value = -14858
$ws.Cells.Item(11, 14).Value = -24284

$ws.Cells.Item(13, 8).Value = 10100.857
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 10100.857
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 10100.857
$ws.Cells.Item(13, 13).ClearContents()
$ws.Cells.Item(13, 14).Value = -10380.857

$ws.Cells.Item(107, 8).Value = 488
$ws.Cells.Item(107, 9).Value = 484
$ws.Cells.Item(107, 10).Value = 500
$ws.Cells.Item(107, 11).Value = 1452
$ws.Cells.Item(107, 12).Value = 1500
$ws.Cells.Item(107, 13).Value = 468
$ws.Cells.Item(107, 14).Value = -5340

$ws.Cells.Item(122, 8).Value = 357969.06
$ws.Cells.Item(122, 9).Value = 385370.53
$ws.Cells.Item(122, 10).Value = 1750
$ws.Cells.Item(122, 11).Value = 1156111.59
$ws.Cells.Item(122, 12).Value = 5250
$ws.Cells.Item(122, 13).Value = -1153661.59
$ws.Cells.Item(122, 14).Value = -10150

$ws.Cells.Item(126, 8).Value = 556516.75
$ws.Cells.Item(126, 9).Value = 769923.9399999999
$ws.Cells.Item(126, 11).Value = 2309771.82
$ws.Cells.Item(126, 13).Value = -2307301.82

$ws.Cells.Item(132, 8).Value = 1268.1464
$ws.Cells.Item(132, 9).Value = 889.6667
$ws.Cells.Item(132, 10).Value = 1665.55
$ws.Cells.Item(132, 11).Value = 2669.0001
$ws.Cells.Item(132, 12).Value = 4996.65
$ws.Cells.Item(132, 13).Value = -139.0001000000002
$ws.Cells.Item(132, 14).Value = -10056.65
